$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.241649031639099
$ws.Range("B1").Value = 1.311608195304871
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.026229858398438
$ws.Range("E1").Value = 0.9293590188026428
